$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9134
$ws1.Range("F4").Value = 473
$ws1.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202402/vF9kexbx1707289709364.jpeg"
$ws1.Range("F5").Value = 455

# Sheet "全部类型" (fourth sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9134
$ws4.Range("F4").Value = 473
$ws4.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202402/vF9kexbx1707289709364.jpeg"
$ws4.Range("F6").Value = 455
